$d = $word.ActiveDocument

# 1. Title replacement - occurs twice (Heading1 and bold byline), both map to the same new text
$d.Content.Find.Execute("Play Bar-X Safecracker Slot for Free - Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Bar-X Safecracker Slot Free Online", 2)

# 2. "What we like" bullet points
$d.Content.Find.Execute("50,000 times the initial bet maximum payout", $true, $false, $false, $false, $false, $true, 1, $false, "Maximum payout of 50,000 times bet", 2)
$d.Content.Find.Execute("Exciting mystery symbols, multipliers, and free spins", $true, $false, $false, $false, $false, $true, 1, $false, "Interesting symbols and their payouts", 2)
$d.Content.Find.Execute("Simple graphics with unique and suspenseful design elements", $true, $false, $false, $false, $false, $true, 1, $false, "Retro graphics and design elements", 2)

# 3. "What we don't like" bullet points
$d.Content.Find.Execute("No progressive jackpot", $true, $false, $false, $false, $false, $true, 1, $false, "No mention of bonus features", 2)
$d.Content.Find.Execute("No interactive bonus round", $true, $false, $false, $false, $false, $true, 1, $false, "Graphics may not appeal to everyone", 2)

# 4. Italic summary line at end of document
$d.Content.Find.Execute("Find out about the features, pros, and cons of Bar-X Safecracker slot game. Play for free and try your perfect heist with this unique retro-themed game.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Bar-X Safecracker Slot and play for free online to win big.", 2)
